$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert new row: 005273382 / MVFC / 100000, right before the 005531186 (RAFAEL) row ---
$ws.Rows(8).Insert()
$ws.Cells.Item(8, 1).NumberFormat = "@"
$ws.Cells.Item(8, 1).Value = "005273382"
$ws.Cells.Item(8, 2).Value = "MVFC"
$ws.Cells.Item(8, 3).Value = 100000

# --- Insert new row: 004387250 / MONICA / 87932.86, right before the 008004995 (JOSE) row ---
$ws.Rows(10).Insert()
$ws.Cells.Item(10, 1).NumberFormat = "@"
$ws.Cells.Item(10, 1).Value = "004387250"
$ws.Cells.Item(10, 2).Value = "MONICA"
$ws.Cells.Item(10, 3).Value = 87932.86

# --- Remove the old 004387250 / MONICA / 58711.51 row (duplicate account, now superseded above) ---
$ws.Rows(15).Delete()

# --- Update account 004466221 (WALTER) -> 004419765, balance 1424.05 -> 1454.63 ---
$ws.Cells.Item(21, 1).NumberFormat = "@"
$ws.Cells.Item(21, 1).Value = "004419765"
$ws.Cells.Item(21, 3).Value = 1454.63

# --- Remove the now-duplicate 004419765 / WALTER / 30.58 row further down the sheet ---
$ws.Rows(184).Delete()
